$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing timing / call-count data points for the "Mumford_1"
# block (rows 21-25), mirroring the existing style used by the analogous
# cells above (D14:D16) -- raw day-fraction time values formatted as h:mm:ss.

$ws.Range("D21").Value = 0.0000115740740740740734993
$ws.Range("D21").NumberFormat = "h:mm:ss"

$ws.Range("D23").Value = 0.0005324074074074074385682
$ws.Range("D23").NumberFormat = "h:mm:ss"
$ws.Range("E23").Value = 9000000

$ws.Range("D24").Value = 0.0002662037037037037192841
$ws.Range("D24").NumberFormat = "h:mm:ss"

$ws.Range("D25").Value = 0.0002777777777777777775368
$ws.Range("D25").NumberFormat = "h:mm:ss"
$ws.Range("E25").Value = 400000

# Move the active selection to reflect where editing finished.
$null = $ws.Range("D26").Select()
